# Auto-generated update of market price / profit columns (H-N) in Kujata_Profits workbook
# Reflects refreshed Universalis market data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3956.625
$ws.Range("I98").Value = 5175
$ws.Range("J98").Value = 301.5
$ws.Range("K98").Value = 5175
$ws.Range("L98").Value = 301.5
$ws.Range("M98").Value = -3677
$ws.Range("N98").Value = -3297.5
$ws.Range("H122").Value = 3956.625
$ws.Range("I122").Value = 5175
$ws.Range("J122").Value = 301.5
$ws.Range("K122").Value = 15525
$ws.Range("L122").Value = 904.5
$ws.Range("M122").Value = -13075
$ws.Range("N122").Value = -5804.5
$ws.Range("H137").Value = 1440.5646
$ws.Range("I137").Value = 1245.5238
$ws.Range("J137").Value = 1850.15
$ws.Range("K137").Value = 3736.5714
$ws.Range("L137").Value = 5550.450000000001
$ws.Range("M137").Value = -1186.5714
$ws.Range("N137").Value = -10650.45
$ws.Range("H138").Value = 1887.7474
$ws.Range("I138").Value = 1068.1428
$ws.Range("J138").Value = 2108.4102
$ws.Range("K138").Value = 3204.4284
$ws.Range("L138").Value = 6325.230599999999
$ws.Range("M138").Value = 1935.5716
$ws.Range("N138").Value = -16605.2306

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2817
$ws.Range("I32").Value = 2817
$ws.Range("K32").Value = 2817
$ws.Range("M32").Value = -2530
$ws.Range("H63").Value = 22224964
$ws.Range("I63").Value = 2318.1316
$ws.Range("K63").Value = 2318.1316
$ws.Range("M63").Value = -1632.1316
$ws.Range("H66").Value = 22224964
$ws.Range("I66").Value = 2318.1316
$ws.Range("K66").Value = 11590.658
$ws.Range("M66").Value = -8158.658000000001
$ws.Range("H74").Value = 2681.1667
$ws.Range("I74").Value = 1978.2307
$ws.Range("K74").Value = 1978.2307
$ws.Range("M74").Value = -1104.2307
$ws.Range("H77").Value = 2681.1667
$ws.Range("I77").Value = 1978.2307
$ws.Range("K77").Value = 9891.1535
$ws.Range("M77").Value = -5523.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1352.2273
$ws.Range("I134").Value = 1167.6471
$ws.Range("K134").Value = 3502.9413
$ws.Range("M134").Value = -967.9412999999995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 220.33333
$ws.Range("I7").Value = 94.59999999999999
$ws.Range("J7").Value = 334.63635
$ws.Range("K7").Value = 94.59999999999999
$ws.Range("L7").Value = 334.63635
$ws.Range("M7").Value = 18.40000000000001
$ws.Range("N7").Value = -560.63635
$ws.Range("H31").Value = 1196.8906
$ws.Range("I31").Value = 1152
$ws.Range("J31").Value = 1630.8334
$ws.Range("K31").Value = 1152
$ws.Range("L31").Value = 1630.8334
$ws.Range("M31").Value = -857
$ws.Range("N31").Value = -2220.8334
$ws.Range("H34").Value = 1196.8906
$ws.Range("I34").Value = 1152
$ws.Range("J34").Value = 1630.8334
$ws.Range("K34").Value = 1152
$ws.Range("L34").Value = 1630.8334
$ws.Range("M34").Value = -950
$ws.Range("N34").Value = -2034.8334
$ws.Range("H99").Value = 1626.8422
$ws.Range("I99").Value = 1536.6666
$ws.Range("K99").Value = 1536.6666
$ws.Range("M99").Value = -38.66660000000002
$ws.Range("H122").Value = 645.75
$ws.Range("I122").Value = 693.05554
$ws.Range("J122").Value = 503.83334
$ws.Range("K122").Value = 2079.16662
$ws.Range("L122").Value = 1511.50002
$ws.Range("M122").Value = 370.83338
$ws.Range("N122").Value = -6411.500019999999
$ws.Range("H126").Value = 1626.8422
$ws.Range("I126").Value = 1536.6666
$ws.Range("K126").Value = 4609.9998
$ws.Range("M126").Value = -2139.9998
$ws.Range("H132").Value = 1551.7646
$ws.Range("I132").Value = 1094.0869
$ws.Range("J132").Value = 2508.7273
$ws.Range("K132").Value = 3282.2607
$ws.Range("L132").Value = 7526.1819
$ws.Range("M132").Value = -752.2606999999998
$ws.Range("N132").Value = -12586.1819
$ws.Range("H134").Value = 18519946
$ws.Range("I134").Value = 1459.7894
$ws.Range("K134").Value = 4379.3682
$ws.Range("M134").Value = -1844.3682

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 14168.8
$ws.Range("I3").Value = 9346
$ws.Range("J3").Value = 18991.6
$ws.Range("K3").Value = 28038
$ws.Range("L3").Value = 56974.8
$ws.Range("M3").Value = -27926
$ws.Range("N3").Value = -57198.8
$ws.Range("H122").Value = 1263.238
$ws.Range("I122").Value = 625.1111
$ws.Range("J122").Value = 1741.8334
$ws.Range("K122").Value = 5625.9999
$ws.Range("L122").Value = 15676.5006
$ws.Range("M122").Value = -3175.9999
$ws.Range("N122").Value = -20576.5006
$ws.Range("H129").Value = 27778704
$ws.Range("I129").Value = 111111470
$ws.Range("J129").Value = 6945513
$ws.Range("K129").Value = 333334410
$ws.Range("L129").Value = 20836539
$ws.Range("M129").Value = -333329410
$ws.Range("N129").Value = -20846539
$ws.Range("H131").Value = 10769569
$ws.Range("I131").Value = 83333750
$ws.Range("J131").Value = 19320.137
$ws.Range("K131").Value = 250001250
$ws.Range("L131").Value = 57960.41099999999
$ws.Range("M131").Value = -249996210
$ws.Range("N131").Value = -68040.41099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 700
$ws.Range("I97").Value = 700
$ws.Range("K97").Value = 700
$ws.Range("M97").Value = -204
$ws.Range("H107").Value = 212.14285
$ws.Range("I107").Value = 215.8
$ws.Range("J107").Value = 203
$ws.Range("K107").Value = 215.8
$ws.Range("L107").Value = 203
$ws.Range("M107").Value = 1704.2
$ws.Range("N107").Value = -4043
$ws.Range("H123").Value = 21614.857
$ws.Range("J123").Value = 21614.857
$ws.Range("L123").Value = 21614.857
$ws.Range("N123").Value = -26514.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 191.5
$ws.Range("I55").Value = 140.35715
$ws.Range("K55").Value = 140.35715
$ws.Range("M55").Value = 32.64285000000001
$ws.Range("H61").Value = 1204
$ws.Range("I61").Value = 1138.1111
$ws.Range("J61").Value = 1322.6
$ws.Range("K61").Value = 1138.1111
$ws.Range("L61").Value = 1322.6
$ws.Range("M61").Value = -936.1111000000001
$ws.Range("N61").Value = -1726.6
$ws.Range("H113").Value = 1204
$ws.Range("I113").Value = 1138.1111
$ws.Range("J113").Value = 1322.6
$ws.Range("K113").Value = 1138.1111
$ws.Range("L113").Value = 1322.6
$ws.Range("M113").Value = 1031.8889
$ws.Range("N113").Value = -5662.6
$ws.Range("H132").Value = 2911.96
$ws.Range("I132").Value = 3289.111
$ws.Range("J132").Value = 2699.8125
$ws.Range("K132").Value = 9867.332999999999
$ws.Range("L132").Value = 8099.4375
$ws.Range("M132").Value = -7337.332999999999
$ws.Range("N132").Value = -13159.4375
$ws.Range("H134").Value = 34920
$ws.Range("J134").Value = 34920
$ws.Range("L134").Value = 34920
$ws.Range("N134").Value = -45060
$ws.Range("H135").Value = 50609.75
$ws.Range("J135").Value = 50609.75
$ws.Range("L135").Value = 50609.75
$ws.Range("N135").Value = -60749.75
$ws.Range("H136").Value = 1585.7059
$ws.Range("I136").Value = 1417.5518
$ws.Range("K136").Value = 4252.6554
$ws.Range("M136").Value = -1702.6554
$ws.Range("H137").Value = 33360
$ws.Range("J137").Value = 33360
$ws.Range("L137").Value = 33360
$ws.Range("N137").Value = -43560
$ws.Range("H138").Value = 35760
$ws.Range("J138").Value = 35760
$ws.Range("L138").Value = 35760
$ws.Range("N138").Value = -46040
$ws.Range("H139").Value = 39260
$ws.Range("J139").Value = 39260
$ws.Range("L139").Value = 39260
$ws.Range("N139").Value = -49540
$ws.Range("H141").Value = 52976
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 52976
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 52976
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -63336

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1748
$ws.Range("I81").Value = 1435.25
$ws.Range("J81").Value = 2999
$ws.Range("K81").Value = 2870.5
$ws.Range("L81").Value = 5998
$ws.Range("M81").Value = -1809.5
$ws.Range("N81").Value = -8120
$ws.Range("H84").Value = 1748
$ws.Range("I84").Value = 1435.25
$ws.Range("J84").Value = 2999
$ws.Range("K84").Value = 14352.5
$ws.Range("L84").Value = 29990
$ws.Range("M84").Value = -9048.5
$ws.Range("N84").Value = -40598
$ws.Range("H113").Value = 393.7619
$ws.Range("I113").Value = 273.375
$ws.Range("J113").Value = 779
$ws.Range("K113").Value = 820.125
$ws.Range("L113").Value = 2337
$ws.Range("M113").Value = 1349.875
$ws.Range("N113").Value = -6677
$ws.Range("H132").Value = 2481.8667
$ws.Range("I132").Value = 2303.3914
$ws.Range("J132").Value = 3068.2856
$ws.Range("K132").Value = 6910.174199999999
$ws.Range("L132").Value = 9204.856800000001
$ws.Range("M132").Value = -4380.174199999999
$ws.Range("N132").Value = -14264.8568
$ws.Range("H135").Value = 92138.336
$ws.Range("J135").Value = 92138.336
$ws.Range("L135").Value = 92138.336
$ws.Range("N135").Value = -102278.336
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("N17").ClearContents()
